# Weekly update: insert a new week's price record for Orégano
# (Mercado Mayorista Lo Valledor de Santiago) as the new first data row,
# pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 309 (shifts old rows 309:365 down to 310:366,
# and grows the sheet dimension from R365 to R366, same as Excel's
# Rows.Insert on a selected row).
$ws.Rows.Item(309).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(309, 1).Value  = 6
$ws.Cells.Item(309, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(309, 3).Value  = "Metropolitana"
$ws.Cells.Item(309, 4).Value  = 45209
$ws.Cells.Item(309, 5).Value  = 13
$ws.Cells.Item(309, 6).Value  = 100112029
$ws.Cells.Item(309, 7).Value  = "Orégano"
$ws.Cells.Item(309, 8).Value  = "Sin especificar"
$ws.Cells.Item(309, 9).Value  = "Primera"
$ws.Cells.Item(309, 10).Value = 33
$ws.Cells.Item(309, 11).Value = 16000
$ws.Cells.Item(309, 12).Value = 16000
$ws.Cells.Item(309, 13).Value = 16000
$ws.Cells.Item(309, 14).Value = "`$/docena de atados"
$ws.Cells.Item(309, 15).Value = "Región Metropolitana"
$ws.Cells.Item(309, 16).Value = 5333
$ws.Cells.Item(309, 17).Value = 3
$ws.Cells.Item(309, 18).Value = "Hortaliza"
